$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.924893617630005
$ws.Range("B1").Value = 1.950632214546204
$ws.Range("C1").Value = 1.858479261398315
$ws.Range("D1").Value = 1.024555087089539
$ws.Range("E1").Value = 0.7011957764625549
